$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D keeps its text representation (values like "1.000", "21.40" etc.
# must not be reinterpreted as numbers with stripped trailing zeros / precision).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.427.42'
$ws.Range("E2").Value = '  +0.84%  '

$ws.Range("D3").Value = '1.880.11'
$ws.Range("E3").Value = '  +1.14%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").Value = '245.09'
$ws.Range("E5").Value = '  +4.98%  '

$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.04%  '

$ws.Range("D7").Value = '0.4764'
$ws.Range("E7").Value = '  +1.88%  '

$ws.Range("D8").Value = '0.2897'
$ws.Range("E8").Value = '  +2.95%  '

$ws.Range("D9").Value = '0.06519'
$ws.Range("E9").Value = '  +0.86%  '

$ws.Range("D10").Value = '21.40'
$ws.Range("E10").Value = '  +1.20%  '

$ws.Range("D11").Value = '0.07721'
$ws.Range("E11").Value = '  +0.36%  '

$ws.Range("D12").Value = '97.19'
$ws.Range("E12").Value = '  +3.89%  '

$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").Value = '0.7412'
$ws.Range("E13").Value = '  +8.73%  '

$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.884.96'
$ws.Range("E14").Value = '  +1.35%  '

$ws.Range("D15").Value = '5.132'
$ws.Range("E15").Value = '  +1.74%  '

$ws.Range("D16").Value = '275.02'
$ws.Range("E16").Value = '  +4.06%  '

$ws.Range("D17").Value = '30.425.34'
$ws.Range("E17").Value = '  +0.98%  '

$ws.Range("D18").Value = '13.54'
$ws.Range("E18").Value = '  +1.03%  '

$ws.Range("D19").Value = '0.000007534'
$ws.Range("E19").Value = '  -1.12%  '

$ws.Range("D20").Value = '1.000'

$ws.Range("D21").Value = '2.122.60'
$ws.Range("E21").Value = '  +0.89%  '

$ws.Range("E22").Value = '  +0.09%  '

$ws.Range("D23").Value = '5.251'
$ws.Range("E23").Value = '  +1.73%  '

$ws.Range("D24").Value = '6.167'
$ws.Range("E24").Value = '  +1.45%  '

$ws.Range("D25").Value = '9.231'
$ws.Range("E25").Value = '  -0.65%  '

$ws.Range("D26").Value = '164.16'
$ws.Range("E26").Value = '  -0.52%  '

$ws.Range("D27").Value = '18.95'
$ws.Range("E27").Value = '  +3.15%  '

$ws.Range("D28").Value = '1.953'
$ws.Range("E28").Value = '  +2.73%  '

$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '1.375'
$ws.Range("E29").Value = '  +0.53%  '

$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").Value = '0.09930'
$ws.Range("E30").Value = '  +1.14%  '

$ws.Range("E31").Value = '  +3.97%  '

$ws.Range("D32").Value = '4.317'
$ws.Range("E32").Value = '  +1.75%  '

$ws.Range("D33").Value = '4.068'
$ws.Range("E33").Value = '  +2.39%  '

$ws.Range("D34").Value = '0.04717'
$ws.Range("E34").Value = '  +1.38%  '

$ws.Range("D35").Value = '1.122'
$ws.Range("E35").Value = '  +0.26%  '

$ws.Range("D36").Value = '0.6971'
$ws.Range("E36").Value = '  +1.17%  '

$ws.Range("D37").Value = '2.719'
$ws.Range("E37").Value = '  +0.23%  '

$ws.Range("D38").Value = '0.01853'
$ws.Range("E38").Value = '  +1.10%  '

$ws.Range("D39").Value = '2.761'
$ws.Range("E39").Value = '  +0.72%  '

$ws.Range("D40").Value = '6.275'
$ws.Range("E40").Value = '  +0.43%  '

$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = '0.4173'
$ws.Range("E41").Value = '  +2.78%  '

$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value = '69.44'
$ws.Range("E42").Value = '  -2.42%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = '1.906'
$ws.Range("E43").Value = '  +0.63%  '

$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D44").Value = '0.9999'
$ws.Range("E44").Value = '  +0.05%  '

$ws.Range("B45").Value = 'TrustWalletToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D45").Value = '0.8390'
$ws.Range("E45").Value = '  +1.31%  '

$ws.Range("D46").Value = '102.11'
$ws.Range("E46").Value = '  -0.14%  '

$ws.Range("D47").Value = '7.091'
$ws.Range("E47").Value = '  +1.99%  '

$ws.Range("D48").Value = '9.194'
$ws.Range("E48").Value = '  +3.13%  '

$ws.Range("D49").Value = '35.24'
$ws.Range("E49").Value = '  +3.36%  '

$ws.Range("D50").Value = '923.15'
$ws.Range("E50").Value = '  -2.91%  '

$ws.Range("D51").Value = '0.05591'
$ws.Range("E51").Value = '  +0.06%  '
